$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B75 was stored as an inline string "3"; convert it to a true numeric value.
$ws.Range("B75").Value = 3

# Append new row 76 with Ying Tang's additional annotation.
$ws.Range("A76").Value = "Ying Tang"

# B76's politeness_score ("4") stays a text value (matches the other text
# columns in the row) rather than being coerced to a number.
$ws.Range("B76").NumberFormat = "@"
$ws.Range("B76").Value = "4"
$ws.Range("B76").Style = "Normal"

$ws.Range("C76").Value = "it would be good"
$ws.Range("D76").Value = "SUG"
$ws.Range("E76").Value = "WRI"
$ws.Range("F76").Value = "85844681-e6c1-4472-a9f5-69a1244b25a4"
$ws.Range("G76").Value = "SktLlGbRZ_annotated.xlsx"
$ws.Range("H76").Value = "Also, it would be good to extend the figure with the second cycle loss."
